$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environment")

# Row 8: "Current velocity" -> "Current speed (m/s)", with tightened limits (1.6 -> 1.0 m/s)
$ws.Range("B8").Value = "Current speed (m/s)"

# Row 9: add code "Hs (m)" in A9 (previously blank), tightened limits (12 -> 5 m)
$ws.Range("A9").Value = "Hs (m)"
$ws.Range("C9").Value = "< 5 m"
$ws.Range("D9").Value = "< 5 m"

$ws.Range("C8").Value = "0.04 - 1.0 m/s"
$ws.Range("D8").Value = "0.04 - 1.0 m/s"

# Rename "Dissolved oxygen concentration at surface" -> "Dissolved oxygen concentration"
$ws.Range("B5").Value = "Dissolved oxygen concentration (mol/m3)"

# Rename "Total chlorophyll concentration at surface" -> "Total chlorophyll concentration"
$ws.Range("B6").Value = "Total chlorophyll concentration (mg/m3)"

# Row 9: "Wave height" -> "Significant wave height (m)"
$ws.Range("B9").Value = "Significant wave height (m)"

$ws.Range("E18").Select()
